$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1986
$ws.Range("J17").Value = 1986
$ws.Range("L17").Value = 5958
$ws.Range("N17").Value = -6294
$ws.Range("H32").Value = 1005.6667
$ws.Range("I32").Value = 999
$ws.Range("K32").Value = 999
$ws.Range("M32").Value = -673
$ws.Range("H98").Value = 2666.3157
$ws.Range("I98").Value = 2893.1177
$ws.Range("J98").Value = 738.5
$ws.Range("K98").Value = 2893.1177
$ws.Range("L98").Value = 738.5
$ws.Range("M98").Value = -1395.1177
$ws.Range("N98").Value = -3734.5
$ws.Range("H100").Value = 619.8
$ws.Range("I100").Value = 462.375
$ws.Range("K100").Value = 462.375
$ws.Range("M100").Value = 78.625
$ws.Range("H103").Value = 2226.9092
$ws.Range("I103").Value = 3142.8333
$ws.Range("J103").Value = 1127.8
$ws.Range("K103").Value = 9428.499899999999
$ws.Range("L103").Value = 3383.4
$ws.Range("M103").Value = -8842.499899999999
$ws.Range("N103").Value = -4555.4
$ws.Range("H112").Value = 2049.5186
$ws.Range("J112").Value = 2089.48
$ws.Range("L112").Value = 6268.440000000001
$ws.Range("N112").Value = -8484.440000000001
$ws.Range("H122").Value = 2666.3157
$ws.Range("I122").Value = 2893.1177
$ws.Range("J122").Value = 738.5
$ws.Range("K122").Value = 8679.3531
$ws.Range("L122").Value = 2215.5
$ws.Range("M122").Value = -6229.3531
$ws.Range("N122").Value = -7115.5
$ws.Range("H137").Value = 1568174.2
$ws.Range("I137").Value = 2778816.2
$ws.Range("K137").Value = 8336448.600000001
$ws.Range("M137").Value = -8333898.600000001
$ws.Range("H138").Value = 3097.7546
$ws.Range("J138").Value = 2812.7354
$ws.Range("L138").Value = 8438.206200000001
$ws.Range("N138").Value = -18718.2062

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 556.7406999999999
$ws.Range("I2").Value = 283.04544
$ws.Range("J2").Value = 1761
$ws.Range("K2").Value = 283.04544
$ws.Range("L2").Value = 1761
$ws.Range("M2").Value = -170.04544
$ws.Range("N2").Value = -1987
$ws.Range("H32").Value = 4430.4443
$ws.Range("I32").Value = 3899.919
$ws.Range("K32").Value = 3899.919
$ws.Range("M32").Value = -3612.919
$ws.Range("H116").Value = 556.7406999999999
$ws.Range("I116").Value = 283.04544
$ws.Range("J116").Value = 1761
$ws.Range("K116").Value = 283.04544
$ws.Range("L116").Value = 1761
$ws.Range("M116").Value = 2010.95456
$ws.Range("N116").Value = -6349
$ws.Range("H123").Value = 86999.60000000001
$ws.Range("J123").Value = 86999.60000000001
$ws.Range("L123").Value = 86999.60000000001
$ws.Range("N123").Value = -96799.60000000001
$ws.Range("H132").Value = 2036.8108
$ws.Range("I132").Value = 1426.2759
$ws.Range("K132").Value = 4278.8277
$ws.Range("M132").Value = -1748.8277

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 556.7406999999999
$ws.Range("I3").Value = 283.04544
$ws.Range("J3").Value = 1761
$ws.Range("K3").Value = 283.04544
$ws.Range("L3").Value = 1761
$ws.Range("M3").Value = -169.04544
$ws.Range("N3").Value = -1989
$ws.Range("H20").Value = 62503692
$ws.Range("I20").Value = 71431940
$ws.Range("K20").Value = 71431940
$ws.Range("M20").Value = -71431693
$ws.Range("H99").Value = 10373.125
$ws.Range("I99").Value = 3938
$ws.Range("K99").Value = 3938
$ws.Range("M99").Value = -2440
$ws.Range("H134").Value = 3984.1843
$ws.Range("I134").Value = 4289.24
$ws.Range("J134").Value = 3397.5386
$ws.Range("K134").Value = 12867.72
$ws.Range("L134").Value = 10192.6158
$ws.Range("M134").Value = -10332.72
$ws.Range("N134").Value = -15262.6158
$ws.Range("H138").Value = 64998.715
$ws.Range("J138").Value = 64998.715
$ws.Range("L138").Value = 64998.715
$ws.Range("N138").Value = -75278.715
$ws.Range("H140").Value = 59016.562
$ws.Range("J140").Value = 59016.562
$ws.Range("L140").Value = 59016.562
$ws.Range("N140").Value = -69376.56200000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3792.8
$ws.Range("I99").Value = 2628.625
$ws.Range("J99").Value = 8449.5
$ws.Range("K99").Value = 2628.625
$ws.Range("L99").Value = 8449.5
$ws.Range("M99").Value = -1130.625
$ws.Range("N99").Value = -11445.5
$ws.Range("H109").Value = 40081
$ws.Range("I109").Value = 21259
$ws.Range("J109").Value = 49492
$ws.Range("K109").Value = 21259
$ws.Range("L109").Value = 49492
$ws.Range("M109").Value = -20219
$ws.Range("N109").Value = -51572
$ws.Range("H126").Value = 3792.8
$ws.Range("I126").Value = 2628.625
$ws.Range("J126").Value = 8449.5
$ws.Range("K126").Value = 7885.875
$ws.Range("L126").Value = 25348.5
$ws.Range("M126").Value = -5415.875
$ws.Range("N126").Value = -30288.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 25962162
$ws.Range("I4").Value = 27393676
$ws.Range("J4").Value = 9738333
$ws.Range("K4").Value = 82181028
$ws.Range("L4").Value = 29214999
$ws.Range("M4").Value = -82180916
$ws.Range("N4").Value = -29215223
$ws.Range("H60").Value = 1114716.5
$ws.Range("J60").Value = 4798.4614
$ws.Range("L60").Value = 14395.3842
$ws.Range("N60").Value = -14897.3842
$ws.Range("H107").Value = 687.6667
$ws.Range("J107").Value = 687.6667
$ws.Range("L107").Value = 2063.0001
$ws.Range("N107").Value = -5903.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H120").Value = 155999
$ws.Range("J120").Value = 155999
$ws.Range("L120").Value = 155999
$ws.Range("N120").Value = -165675
$ws.Range("H126").Value = 2500
$ws.Range("I126").Value = 2500
$ws.Range("K126").Value = 7500
$ws.Range("M126").Value = -5030

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2688.4211
$ws.Range("J7").Value = 1834.8572
$ws.Range("L7").Value = 1834.8572
$ws.Range("N7").Value = -2058.8572
$ws.Range("H93").Value = 929.1177
$ws.Range("I93").Value = 907.3077
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 907.3077
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 340.6923
$ws.Range("N93").Value = -3496
$ws.Range("H126").Value = 2688.4211
$ws.Range("J126").Value = 1834.8572
$ws.Range("L126").Value = 5504.571599999999
$ws.Range("N126").Value = -10444.5716
$ws.Range("H136").Value = 4540.6665
$ws.Range("I136").Value = 4811.125
$ws.Range("K136").Value = 14433.375
$ws.Range("M136").Value = -11883.375

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H132").Value = 2354.7026
$ws.Range("I132").Value = 2166.0293
$ws.Range("K132").Value = 6498.0879
$ws.Range("M132").Value = -3968.0879
$ws.Range("N128").ClearContents()
